$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "AdaBoostClassifier"
$ws.Range("B10").Value = 0.9947049924357034
$ws.Range("C10").Value = 0.9947424529098876
$ws.Range("D10").Value = 0.9947026872498572
$ws.Range("E10").Value = 0.9947039454146909
$ws.Range("F10").Value = 0.4025695323944092
$ws.Range("G10").Value = 0.01609516143798828
$ws.Range("H10").Value = 0.00001217485736610309
$ws.Range("I10").Value = 0.9927355278093076
$ws.Range("J10").Value = 0.005448354143019296
